# Integrate read_batch and simulate_synthetic functions:
# add an English-translation column (C) next to the existing Russian
# condition-label column (B), matching each row's translated text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Russian -> English translation lookup for the condition labels used in column B.
$translations = @{
    "Покой"                           = "Rest";
    "постоянное легкое моргание"      = "Continuous light blinking";
    "постоянное сильное моргание"     = "Continuous strong blinking";
    "постоянное сжимание челюсти"     = "Continuous jaw clenching";
    "покой с морганием слабым"        = "Rest with weak blinking";
    "покой с морганием сильным"       = "Rest with strong blinking";
    "покой с сжиманием челюсти"       = "Rest with jaw clenching";
    "покой со всем вместе"            = "Rest with all together";
}

# Data rows 2-26 (row 22 is an intentional blank separator row, matching
# the gap already present between the "постоянное..." block and the
# "покой с..." block).
for ($row = 2; $row -le 26; $row++) {
    $label = $ws.Cells.Item($row, 2).Value2
    if ($label -eq $null -or $label -eq "") {
        continue
    }
    $english = $translations[$label]
    if ($english -eq $null) {
        continue
    }

    $ws.Cells.Item($row, 3).Value2 = $english

    # Match the formatting/style already used on column A for this row
    # (e.g. the "s=1" cell style) by copying formats across.
    $ws.Cells.Item($row, 1).Copy() | Out-Null
    $ws.Cells.Item($row, 3).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
}

$excel.CutCopyMode = $false

# New column C needs its own width, similar in spirit to column B's custom width.
$ws.Columns.Item(3).ColumnWidth = 29.3

# The two pictures on the sheet swap their display names.
$ws.Shapes.Item(1).Name = "image1.png"
$ws.Shapes.Item(2).Name = "image2.png"
